$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.921.69"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "3.384.39"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "571.94"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "136.78"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.382.90"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "0.469"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").Value = "7.63"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "3.960.15"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "26.40"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "3.381.31"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("D18").Value = "60.971.33"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "13.85"
$ws.Range("E19").Value = "  -1.66%  "
$ws.Range("D20").Value = "5.83"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "9.28"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "374.69"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "3.508.75"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "0.549"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "70.70"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -2.59%  "
$ws.Range("D28").Value = "1.61"
$ws.Range("E28").Value = "  -6.64%  "
$ws.Range("E29").Value = "  +6.83%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").Value = "7.35"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("D33").Value = "2.13"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "23.35"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "5.11"
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "6.80"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "164.88"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").Value = "0.0765"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").Value = "25.70"
$ws.Range("E41").Value = "  +4.81%  "
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "0.773"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "41.92"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("D46").Value = "4.36"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "1.17"
$ws.Range("E47").Value = "  -4.25%  "
$ws.Range("D48").Value = "2.514.41"
$ws.Range("E48").Value = "  +7.09%  "
$ws.Range("D49").Value = "23.52"
$ws.Range("E49").Value = "  +3.55%  "
$ws.Range("D50").Value = "6.76"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  +1.83%  "
